$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain literal values are assigned directly. Cells whose new text looks
# like a plain number (e.g. '596.95') are prefixed with a leading apostrophe
# so Excel keeps storing them as text (matching the original inlineStr
# cells in this sheet) instead of auto-converting them to numbers.

$ws.Range("D2").Value = "71.939.00"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "2.681.79"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'596.95"
$ws.Range("E5").Value = "  -2.21%  "
$ws.Range("D6").Value = "'173.84"
$ws.Range("E6").Value = "  -4.25%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").Value = "2.682.96"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("E10").Value = "  -4.58%  "
$ws.Range("E11").Value = "  +2.03%  "
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("D13").Value = "'4.99"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").Value = "3.175.97"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "71.957.53"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000184"
$ws.Range("E16").Value = "  -3.13%  "
$ws.Range("D17").Value = "'26.13"
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("D18").Value = "2.681.96"
$ws.Range("E18").Value = "  +1.69%  "
$ws.Range("D19").Value = "'12.22"
$ws.Range("E19").Value = "  +5.43%  "
$ws.Range("D20").Value = "'8.19"
$ws.Range("E20").Value = "  +3.96%  "
$ws.Range("D21").Value = "'370.92"
$ws.Range("E21").Value = "  -3.94%  "
$ws.Range("D22").Value = "'4.19"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "'1.99"
$ws.Range("E23").Value = "  -2.56%  "
$ws.Range("D24").Value = "'72.27"
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "'4.33"
$ws.Range("E26").Value = "  -2.99%  "
$ws.Range("D27").Value = "'9.74"
$ws.Range("E27").Value = "  -2.66%  "
$ws.Range("D28").Value = "2.821.44"
$ws.Range("E28").Value = "  +1.48%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "0.0₃0965"
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").Value = "'8.05"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").Value = "'499.23"
$ws.Range("E32").Value = "  -9.80%  "
$ws.Range("D33").Value = "'1.29"
$ws.Range("E33").Value = "  -3.65%  "
$ws.Range("E34").Value = "  -1.60%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'164.48"
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("D37").Value = "'19.62"
$ws.Range("E37").Value = "  +1.51%  "
$ws.Range("D38").Value = "'19.09"
$ws.Range("D39").Value = "'1.37"
$ws.Range("E39").Value = "  -2.55%  "
$ws.Range("E40").Value = "  -7.18%  "
$ws.Range("D41").Value = "'1.77"
$ws.Range("E41").Value = "  -5.43%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "'5.00"
$ws.Range("E43").Value = "  -2.25%  "
$ws.Range("D44").Value = "'0.333"
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("E45").Value = "  -3.68%  "
$ws.Range("D46").Value = "'157.31"
$ws.Range("E46").Value = "  +3.55%  "
$ws.Range("D47").Value = "'39.37"
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("D48").Value = "'0.564"
$ws.Range("E48").Value = "  +3.95%  "
$ws.Range("D49").Value = "'3.73"
$ws.Range("E49").Value = "  +1.13%  "
$ws.Range("D50").Value = "'1.75"
$ws.Range("E50").Value = "  +2.57%  "
$ws.Range("D51").Value = "'0.0759"
$ws.Range("E51").Value = "  +0.04%  "
